$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("S30", "S31")
Write-Host "done"
